# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# Forcing Text number format so numeric-looking strings keep their exact
# textual representation (matching the source inlineStr cells), then
# restoring the default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "35.165.70"
Set-TextValue "E2" "  -0.13%  "
Set-TextValue "D3" "1.904.67"
Set-TextValue "E3" "  +0.48%  "
Set-TextValue "E4" "  -0.44%  "
Set-TextValue "D5" "253.84"
Set-TextValue "E5" "  +3.30%  "
Set-TextValue "D6" "0.699"
Set-TextValue "E6" "  +2.26%  "
Set-TextValue "E7" "  -0.39%  "
Set-TextValue "D8" "41.73"
Set-TextValue "E8" "  +2.38%  "
Set-TextValue "D9" "0.355"
Set-TextValue "E9" "  +2.36%  "
Set-TextValue "D10" "52.69"
Set-TextValue "E10" "  +0.66%  "
Set-TextValue "D11" "0.0759"
Set-TextValue "E11" "  +5.53%  "
Set-TextValue "D12" "0.0979"
Set-TextValue "E12" "  -0.29%  "
Set-TextValue "D13" "13.19"
Set-TextValue "E13" "  +3.20%  "
Set-TextValue "D14" "2.181.60"
Set-TextValue "E14" "  +0.51%  "
Set-TextValue "D15" "0.737"
Set-TextValue "E15" "  +4.59%  "
Set-TextValue "D16" "5.03"
Set-TextValue "E16" "  +4.91%  "
Set-TextValue "D17" "1.898.21"
Set-TextValue "E17" "  +0.25%  "
Set-TextValue "D18" "35.152.97"
Set-TextValue "E18" "  -0.17%  "
Set-TextValue "D19" "73.90"
Set-TextValue "E19" "  +2.73%  "
Set-TextValue "E20" "  +3.25%  "
Set-TextValue "D21" "243.19"
Set-TextValue "D22" "13.09"
Set-TextValue "E22" "  +3.46%  "
Set-TextValue "D23" "5.06"
Set-TextValue "E23" "  +5.73%  "
Set-TextValue "E24" "  -0.37%  "
Set-TextValue "E25" "  +5.79%  "
Set-TextValue "D26" "2.33"
Set-TextValue "E26" "  +1.01%  "
Set-TextValue "D27" "167.68"
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "8.61"
Set-TextValue "E28" "  +0.10%  "
Set-TextValue "D29" "18.55"
Set-TextValue "E29" "  +1.31%  "
Set-TextValue "E30" "  +0.22%  "
Set-TextValue "D31" "4.128.39"
Set-TextValue "E31" "  -0.34%  "
Set-TextValue "E32" "  +6.34%  "
Set-TextValue "D33" "0.0601"
Set-TextValue "E33" "  +6.07%  "
Set-TextValue "D34" "4.35"
Set-TextValue "E34" "  +4.85%  "
Set-TextValue "E35" "  +8.95%  "
Set-TextValue "D36" "4.26"
Set-TextValue "E36" "  +3.59%  "
Set-TextValue "E37" "  -0.46%  "
Set-TextValue "D38" "0.855"
Set-TextValue "E38" "  -6.78%  "
Set-TextValue "E39" "  -0.16%  "
Set-TextValue "D40" "99.04"
Set-TextValue "E40" "  +10.07%  "
Set-TextValue "D41" "17.14"
Set-TextValue "E41" "  +4.43%  "
Set-TextValue "E42" "  +4.45%  "
Set-TextValue "E43" "  +1.98%  "
Set-TextValue "D44" "0.0656"
Set-TextValue "E44" "  +3.09%  "
Set-TextValue "E45" "  +0.17%  "
Set-TextValue "D46" "1.309.78"
Set-TextValue "E46" "  -2.73%  "
Set-TextValue "E47" "  +0.11%  "
Set-TextValue "E48" "  -1.03%  "
Set-TextValue "E49" "  +1.97%  "
Set-TextValue "E50" "  +2.27%  "
Set-TextValue "D51" "0.0754"
Set-TextValue "E51" "  +7.38%  "
